$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet 3: Inflation rate (%) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Inflation rate (%)"

$ws3.Range("A1").Value = "date"
$ws3.Range("B1").Value = "Inflation rate (%)"

$ws3.Range("A2").Value = 2015
$ws3.Range("B2").Value = 15.5
$ws3.Range("A3").Value = 2016
$ws3.Range("B3").Value = 7
$ws3.Range("A4").Value = 2017
$ws3.Range("B4").Value = 3.7
$ws3.Range("A5").Value = 2018
$ws3.Range("B5").Value = 2.9
$ws3.Range("A6").Value = 2019
$ws3.Range("B6").Value = 4.5

$ws.Range("A1:B1").Copy()
$ws3.Range("A1:B1").PasteSpecial(-4122)

# --- Sheet 4: Current account balance (% of GDP) ---
# Sheet (tab) names are capped at 31 chars in Excel, so the tab name is
# truncated even though the header cell/shared string keeps the full text.
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Current account balance (% of G"

$ws4.Range("A1").Value = "date"
$ws4.Range("B1").Value = "Current account balance (% of GDP)"

$ws4.Range("A2").Value = 2015
$ws4.Range("B2").Value = 5
$ws4.Range("A3").Value = 2016
$ws4.Range("B3").Value = 1.9
$ws4.Range("A4").Value = 2017
$ws4.Range("B4").Value = 2
$ws4.Range("A5").Value = 2018
$ws4.Range("B5").Value = 6.9
$ws4.Range("A6").Value = 2019
$ws4.Range("B6").Value = 3.8

$ws.Range("A1:B1").Copy()
$ws4.Range("A1:B1").PasteSpecial(-4122)

# Restore original active sheet/selection
$ws.Activate()
